# Applies the cryptos-list price/volume refresh described by the commit diff.
# Column D ("Price") holds plain-text numbers (some with thousands separators like
# "67.876.06" that are not valid numeric literals) and must stay text, so numeric-
# looking D values are entered with a leading apostrophe (forces text entry, same as
# typing it in Excel) and the style is reset to "Normal" right after so no stray
# quote-prefix / text-format style lingers on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '67.876.06'
$ws.Range("E2").Value = '  -2.16%  '
# Row 3: Ethereum
$ws.Range("D3").Value = '3.802.22'
$ws.Range("E3").Value = '  +0.82%  '
# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.08%  '
# Row 5: BNB
$ws.Range("D5").Value = "'601.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.12%  '
# Row 6: Solana
$ws.Range("D6").Value = "'171.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.69%  '
# Row 7: LidoStakedEther
$ws.Range("D7").Value = '3.801.32'
$ws.Range("E7").Value = '  +0.88%  '
# Row 8: USDC
$ws.Range("E8").Value = '  +0.06%  '
# Row 9: XRP
$ws.Range("E9").Value = '  +0.35%  '
# Row 10: Dogecoin
$ws.Range("D10").Value = "'0.159"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.58%  '
# Row 11: Toncoin
$ws.Range("D11").Value = "'6.19"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.17%  '
# Row 12: Cardano
$ws.Range("E12").Value = '  -3.77%  '
# Row 13: Avalanche
$ws.Range("D13").Value = "'38.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.99%  '
# Row 14: ShibaInu
$ws.Range("D14").Value = "'0.0000244"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.80%  '
# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '4.442.42'
$ws.Range("E15").Value = '  +1.07%  '
# Row 16: WrappedEther
$ws.Range("D16").Value = '3.804.99'
# Row 17: WrappedBTC
$ws.Range("D17").Value = '67.918.80'
# Row 18: Polkadot
$ws.Range("D18").Value = "'7.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.08%  '
# Row 19: TRON
$ws.Range("E19").Value = '  -3.88%  '
# Row 20: Chainlink
$ws.Range("E20").Value = '  +5.63%  '
# Row 21: BitcoinCash
$ws.Range("D21").Value = "'494.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.89%  '
# Row 22: Uniswap
$ws.Range("E22").Value = '  -1.52%  '
# Row 23: Polygon
$ws.Range("D23").Value = "'0.741"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.39%  '
# Row 24: Litecoin
$ws.Range("D24").Value = "'85.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.95%  '
# Row 25: Fetch.AI
$ws.Range("E25").Value = '  -4.73%  '
# Row 26: PEPE
$ws.Range("E26").Value = '  +6.21%  '
# Row 27: InternetComputer(DFINITY)
$ws.Range("E27").Value = '  -3.60%  '
# Row 28: RenderToken
$ws.Range("E28").Value = '  -3.91%  '
# Row 29: Dai
$ws.Range("E29").Value = '  +0.09%  '
# Row 30: PancakeSwap
$ws.Range("E30").Value = '  -0.44%  '
# Row 31: ImmutableX
$ws.Range("E31").Value = '  -3.24%  '
# Row 32: EthereumClassic
$ws.Range("D32").Value = "'33.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.87%  '
# Row 33: NEARProtocol
$ws.Range("E33").Value = '  -2.07%  '
# Row 34: Hedera
$ws.Range("E34").Value = '  -3.68%  '
# Row 35: FirstDigitalUSD
$ws.Range("E35").Value = '  +0.05%  '
# Row 36: Mantle
$ws.Range("E36").Value = '  -3.05%  '
# Row 37: Filecoin
$ws.Range("E37").Value = '  -4.72%  '
# Row 38: Bittensor
$ws.Range("B38").Value = 'Bittensor'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D38").Value = "'463.57"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.74%  '
# Row 39: TheGraph
$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").Value = "'0.332"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.73%  '
# Row 40: Kaspa
$ws.Range("D40").Value = "'0.132"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.24%  '
# Row 41: Stacks
$ws.Range("E41").Value = '  -2.79%  '
# Row 43: dogwifhat
$ws.Range("E43").Value = '  -5.07%  '
# Row 44: Cosmos
$ws.Range("D44").Value = "'8.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.59%  '
# Row 45: USDe
$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").Value = "'1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.01%  '
# Row 46: Arweave
$ws.Range("B46").Value = 'Arweave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D46").Value = "'40.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -9.03%  '
# Row 47: Maker
$ws.Range("D47").Value = '2.844.88'
$ws.Range("E47").Value = '  -3.78%  '
# Row 48: EnergySwap
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = "'25.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +19.18%  '
# Row 49: Monero
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").Value = "'139.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.41%  '
# Row 50: VeChain
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").Value = "'0.0352"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.42%  '
# Row 51: InjectiveProtocol
$ws.Range("D51").Value = "'25.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.29%  '
